# Flora_009 commit: "Add unipa transcriptions (-Hasan 026)"
#
# Real content changes behind the noisy re-save diff:
#   1. Sheet "Sampling Events" (sheet1) cell B2: fix a stray space in the
#      transcribed UNIPA catalogue code
#        "UNIPA- 2005DR-AM009-PM001" -> "UNIPA-2005DR-AM009-PM001"
#   2. Sheet "Occurrences" (sheet2) column A (rows 2-4): same catalogue-code
#      typo fix.
#   3. Sheet "Occurrences" (sheet2): drop the large block of left-over,
#      empty-but-formatted rows (5-56) so the sheet's used range shrinks
#      back down to A1:Q4.

$wb = $excel.ActiveWorkbook

$samplingEvents = $wb.Worksheets.Item("Sampling Events")
$occurrences    = $wb.Worksheets.Item("Occurrences")

# 1. Fix the transcription typo on "Sampling Events"
$samplingEvents.Range("B2").Value = "UNIPA-2005DR-AM009-PM001"

# 2. Fix the same typo on "Occurrences" (rows 2, 3, 4 all share this value)
$occurrences.Range("A2:A4").Value = "UNIPA-2005DR-AM009-PM001"

# 3. Remove the stray empty formatted rows (5-56) trailing the real data
$occurrences.Range("A5:A56").EntireRow.Delete()
